$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.434.10"
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.502.57"
$ws.Range("E3").Value = "  -5.44%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.12"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.92"
$ws.Range("E6").Value = "  -5.19%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.500.91"
$ws.Range("E9").Value = "  -5.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -7.49%  "
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  -4.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.83"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.960.89"
$ws.Range("E14").Value = "  -5.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.343.80"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  -6.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.68"
$ws.Range("E17").Value = "  -4.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.514.76"
$ws.Range("E18").Value = "  -5.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.34"
$ws.Range("E19").Value = "  -7.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.71"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "347.64"
$ws.Range("E21").Value = "  -6.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.92"
$ws.Range("E22").Value = "  -5.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.95"
$ws.Range("E23").Value = "  -5.88%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.48"
$ws.Range("E25").Value = "  -3.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.98"
$ws.Range("E26").Value = "  -7.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.93"
$ws.Range("E27").Value = "  -7.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.633.81"
$ws.Range("E28").Value = "  -5.28%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0894"
$ws.Range("E30").Value = "  -7.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "470.30"
$ws.Range("E32").Value = "  -6.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.26"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.74"
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.115"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.65"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.93"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.41"
$ws.Range("E39").Value = "  -4.82%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.71"
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.317"
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.60"
$ws.Range("E43").Value = "  -8.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.17"
$ws.Range("E44").Value = "  -13.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("E45").Value = "  -10.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.12"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.85"
$ws.Range("E47").Value = "  -6.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.526"
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.51"
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.59"
$ws.Range("E50").Value = "  -6.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0730"
$ws.Range("E51").Value = "  -2.63%  "
